$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update English-column (C) translations that changed wording --------
# (rows where the English text is unchanged are left untouched)
$ws.Range("C2").Value2  = "Repository is null"
$ws.Range("C6").Value2  = "Repository"
$ws.Range("C8").Value2  = "Tag"
$ws.Range("C12").Value2 = "Delete the repository"
$ws.Range("C13").Value2 = "Confirm to delete the repository"
$ws.Range("C15").Value2 = "The repository cannot be null. Please manually clear the image or check the check box to delete images in the repository by force'"
$ws.Range("C16").Value2 = "Repository"

# --- New column D: reserved/empty column with its own (red font) style --
$colD = $ws.Range("D1:D17")
$colD.Font.Name = "Times New Roman"
$colD.Font.Size = 12
$colD.Font.Color = 255
$colD.WrapText = $true

# --- Column widths ---------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 17.375
$ws.Columns.Item(2).ColumnWidth = 48.125
$ws.Columns.Item(3).ColumnWidth = 48.125
$ws.Columns.Item(4).ColumnWidth = 48.125

# --- Page setup (paper size / orientation) ---------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection moved to C19 --------------------------------------------
$ws.Range("C19").Select()
